$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.014.52"
$ws.Range("E2").Value = "  -5.64%  "
$ws.Range("D3").Value = "3.098.42"
$ws.Range("E3").Value = "  -6.60%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.09"
$ws.Range("E5").Value = "  -7.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.31"
$ws.Range("E6").Value = "  -6.76%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.100.19"
$ws.Range("E8").Value = "  -6.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.443"
$ws.Range("E9").Value = "  -6.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("E10").Value = "  -9.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  -10.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -7.99%  "
$ws.Range("D13").Value = "3.643.99"
$ws.Range("E13").Value = "  -6.17%  "
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.21"
$ws.Range("E15").Value = "  -6.95%  "
$ws.Range("D16").Value = "3.105.76"
$ws.Range("E16").Value = "  -6.24%  "
$ws.Range("D17").Value = "57.065.93"
$ws.Range("E17").Value = "  -5.50%  "
$ws.Range("E18").Value = "  -10.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.71"
$ws.Range("E19").Value = "  -7.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -11.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.90"
$ws.Range("E21").Value = "  -8.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "341.14"
$ws.Range("E22").Value = "  -9.02%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.11"
$ws.Range("E24").Value = "  -8.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.501"
$ws.Range("E25").Value = "  -8.14%  "
$ws.Range("D26").Value = "3.243.02"
$ws.Range("E26").Value = "  -5.63%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.165"
$ws.Range("E28").Value = "  -4.77%  "
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("E29").Value = "  -12.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.66"
$ws.Range("E31").Value = "  -8.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.92"
$ws.Range("E32").Value = "  -10.12%  "
$ws.Range("E33").Value = "  -9.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.33"
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("E35").Value = "  -5.80%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.45"
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.73"
$ws.Range("E37").Value = "  -9.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.11"
$ws.Range("E38").Value = "  -9.76%  "
$ws.Range("E39").Value = "  -11.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.11"
$ws.Range("E40").Value = "  -6.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0681"
$ws.Range("E41").Value = "  -8.28%  "
$ws.Range("D42").Value = "3.131.22"
$ws.Range("E42").Value = "  -6.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.19"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.676"
$ws.Range("E44").Value = "  -10.30%  "
$ws.Range("E45").Value = "  -7.89%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").Value = "  -7.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.44"
$ws.Range("E48").Value = "  -9.86%  "
$ws.Range("D49").Value = "2.248.88"
$ws.Range("E49").Value = "  -4.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.12"
$ws.Range("E50").Value = "  -6.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.68"
$ws.Range("E51").Value = "  -8.71%  "
